# Apply updated values to result_data_KNN.xlsx (terrestrial_mammals, combination_1_ABCD, AC/15/seed5)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C11"  = -12.618
    "A12"  = -21.544
    "C23"  = -12.309
    "A27"  = -21.8
    "C28"  = -12.992
    "A32"  = -22.023
    "C32"  = -13.538
    "C34"  = -12.345
    "A36"  = -20.43
    "A38"  = -19.83
    "C42"  = -12.36
    "A46"  = -21.791
    "C49"  = -13.271
    "A54"  = -21.945
    "C54"  = -12.803
    "A55"  = -22.173
    "A56"  = -21.928
    "A67"  = -21.577
    "A69"  = -21.503
    "A72"  = -21.689
    "C78"  = -12.56
    "C80"  = -12.664
    "A83"  = -21.805
    "A86"  = -22.076
    "A91"  = -20.755
    "A93"  = -21.452
    "C97"  = -11.516
    "A99"  = -21.708
    "C99"  = -12.418
    "C101" = -12.188
    "A104" = -21.175
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
